$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.312.84'
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').Value = '1.834.36'
$ws.Range('E4').Value = '  +0.90%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.009'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3689'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8863'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.45'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.17%  '
$ws.Range('D12').Value = '1.881.08'
$ws.Range('E12').Value = '  +4.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07330'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.436'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.12%  '
$ws.Range('E15').Value = '  +2.96%  '
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008793'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.81'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.23%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '27.526.27'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.287'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('D24').Value = '2.090.57'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.897'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('E27').Value = '  +1.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.153'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.227'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08990'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.99%  '
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.175'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.549'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.950'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.011'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.100'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05344'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01957'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.967'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.394'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.232'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5310'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1659'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.483'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4934'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.90%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.672'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06298'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.09%  '
